# ENH: Add dynamic lapse
# Applies the "dynamic lapse" enhancement to model_parameters.xlsx:
#  - ParamList:  fix two descriptions' wording, add new "is_lapse_dynamic" parameter row
#  - SpaceParams: add a new "is_lapse_dynamic" column (TRUE for all products)
#  - GMXB: add "dyn_lapse_param_id" and "dyn_lapse_floor" columns with per-row ids/floors

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ParamList
# ---------------------------------------------------------------------------
$paramList = $wb.Worksheets.Item("ParamList")

# Wording tweaks ("str." -> "string")
$paramList.Range("C11").Value = "Sensitivity ID in string"
$paramList.Range("C12").Value = "Interest rate sensitivity ID in string"

# New parameter row describing the dynamic-lapse switch
$paramList.Range("A18").Value = "is_lapse_dynamic"
$paramList.Range("B18").Value = "SPACE"
$paramList.Range("C18").Value = "Wheter the lapse assumption is dynamic in boolean"

# Column C widened to fit the new/longer text (closest reachable width to 43.53515625)
$paramList.Columns.Item(3).ColumnWidth = 42.65

# ---------------------------------------------------------------------------
# SpaceParams
# ---------------------------------------------------------------------------
$spaceParams = $wb.Worksheets.Item("SpaceParams")

$spaceParams.Range("E1").Value = "is_lapse_dynamic"
$spaceParams.Range("E2").Value = $true
$spaceParams.Range("E3").Value = $true
$spaceParams.Range("E4").Value = $true

# New column width (closest reachable width to 17.84375)
$spaceParams.Columns.Item(5).ColumnWidth = 17.0

# ---------------------------------------------------------------------------
# GMXB
# ---------------------------------------------------------------------------
$gmxb = $wb.Worksheets.Item("GMXB")

$gmxb.Range("P1").Value = "dyn_lapse_param_id"
$gmxb.Range("Q1").Value = "dyn_lapse_floor"

$gmxb.Range("P2").Value = "DL001A"
$gmxb.Range("Q2").Value = 0
$gmxb.Range("Q2").NumberFormat = "0%"

$gmxb.Range("P3").Value = "DL001B"
$gmxb.Range("Q3").Value = 0
$gmxb.Range("Q3").NumberFormat = "0%"

$gmxb.Range("P4").Value = "DL002A"
$gmxb.Range("Q4").Value = 0.03
$gmxb.Range("Q4").NumberFormat = "0%"

$gmxb.Range("P5").Value = "DL002B"
$gmxb.Range("Q5").Value = 0.05
$gmxb.Range("Q5").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# Selections (match the cells the author had selected when saving)
# ---------------------------------------------------------------------------
$runParams = $wb.Worksheets.Item("RunParams")
$runParams.Activate()
$runParams.Range("I9").Select()

$gmxb.Activate()
$gmxb.Range("Q1").Select()

$spaceParams.Activate()
$spaceParams.Range("E1").Select()

# ParamList stays the active/visible tab, selected last
$paramList.Activate()
$paramList.Range("C12").Select()
